$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.261385560035706
$ws.Range("B1").Value = 2.487245559692383
$ws.Range("C1").Value = 4.66655969619751
$ws.Range("D1").Value = 2.016480922698975
$ws.Range("E1").Value = 1.149526596069336
